$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07477033333333334
$ws.Range("H2").Value = 0.224311
$ws.Range("M2").Value = 2.680851666666667
$ws.Range("N2").Value = 8.042555
$ws.Range("O2").Value = 0.1074910720871699
$ws.Range("P2").Value = 0.1074910720871699
$ws.Range("Q2").Value = 0.2004481727338889
$ws.Range("R2").Value = 1.804033554605
$ws.Range("S2").Value = 0.1074910720871699
$ws.Range("T2").Value = 0.1074910720871699

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07477033333333334
$ws.Range("H3").Value = 0.224311
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("O3").Value = 0.7423457357290222
$ws.Range("P3").Value = 0.7423457357290222
$ws.Range("Q3").Value = 1.384318189170223
$ws.Range("R3").Value = 12.458863702532
$ws.Range("S3").Value = 0.7423457357290222
$ws.Range("T3").Value = 0.7423457357290222

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07477033333333334
$ws.Range("H4").Value = 0.224311
$ws.Range("M4").Value = 3.745104
$ws.Range("N4").Value = 11.235312
$ws.Range("O4").Value = 0.1501631921838079
$ws.Range("P4").Value = 0.1501631921838079
$ws.Range("Q4").Value = 0.2800226744480001
$ws.Range("R4").Value = 2.520204070032
$ws.Range("S4").Value = 0.1501631921838079
$ws.Range("T4").Value = 0.1501631921838079
